$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-16, 0, 1),
    @(-11.314, -11.314, 1),
    @(11.314, -11.314, 1),
    @(0, -16, 1),
    @(-11.314, 11.314, 1),
    @(11.314, 11.314, 1),
    @(0, 16, 1),
    @(0, -16, 1),
    @(-16, 0, 1),
    @(-11.314, -11.314, 1)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("D19").Select()
